# Automatic update of files.
# Update row 2, row 3 and row 4 of the "Artfynd" sheet with revised
# species-occurrence data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: was Garnlav / Alectoria sarmentosa -> becomes Lunglav / Lobaria pulmonaria
$ws.Range("A2").Value = 112252574
$ws.Range("B2").Value = 78699
$ws.Range("E2").Value = 6458
$ws.Range("F2").Value = "Lunglav"
$ws.Range("G2").Value = "Lobaria pulmonaria"
$ws.Range("H2").Value = "(L.) Hoffm."
$ws.Range("AC2").ClearContents()

# --- Row 3: was Rosenticka / Rhodofomes roseus -> becomes Garnlav / Alectoria sarmentosa
$ws.Range("A3").Value = 112252573
$ws.Range("B3").Value = 77636
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("AC3").Value = "Ganska rikligt"

# --- Row 4: was Lunglav / Lobaria pulmonaria -> becomes Rosenticka / Rhodofomes roseus
$ws.Range("A4").Value = 112252575
$ws.Range("B4").Value = 89820
$ws.Range("E4").Value = 658
$ws.Range("F4").Value = "Rosenticka"
$ws.Range("G4").Value = "Rhodofomes roseus"
$ws.Range("H4").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("AC4").Value = "På två granlågor, gyttrad på en."
